$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 41: KKR vs DC (row 53)
$ws.Range("E53").Value = 60
$ws.Range("H53").Value = 80
$ws.Range("K53").Value = 100
$ws.Range("N53").Value = 30
$ws.Range("Q53").Value = 70
$ws.Range("T53").Value = 20
$ws.Range("W53").Value = 50
$ws.Range("Z53").Value = 0
$ws.Range("AC53").Value = 40

# Contest 42: MI vs PBKS (row 54)
$ws.Range("E54").Value = 100
$ws.Range("H54").Value = 70
$ws.Range("K54").Value = 60
$ws.Range("N54").Value = 40
$ws.Range("Q54").Value = 80
$ws.Range("T54").Value = 30
$ws.Range("W54").Value = 0
$ws.Range("Z54").Value = 50
$ws.Range("AC54").Value = 20

$excel.CalculateFullRebuild()
$wb.Save()
